$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3479.2222
$ws.Range("I113").Value = 3273.2856
$ws.Range("K113").Value = 3273.2856
$ws.Range("M113").Value = -19.28560000000016

$ws.Range("H137").Value = 57360228
$ws.Range("J137").Value = 2952356.8
$ws.Range("L137").Value = 8857070.399999999
$ws.Range("N137").Value = -8862170.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2269.2307
$ws.Range("I45").Value = 2275
$ws.Range("K45").Value = 2275
$ws.Range("M45").Value = -1898

$ws.Range("H74").Value = 2978693.8
$ws.Range("I74").Value = 3677556.5
$ws.Range("K74").Value = 3677556.5
$ws.Range("M74").Value = -3676682.5

$ws.Range("H77").Value = 2978693.8
$ws.Range("I77").Value = 3677556.5
$ws.Range("K77").Value = 18387782.5
$ws.Range("M77").Value = -18383414.5

$ws.Range("H102").Value = 31900.584
$ws.Range("I102").Value = 34394
$ws.Range("K102").Value = 34394
$ws.Range("M102").Value = -32772

$ws.Range("H122").Value = 3885.8635
$ws.Range("I122").Value = 3838.3333
$ws.Range("K122").Value = 11514.9999
$ws.Range("M122").Value = -9064.999899999999

$ws.Range("H132").Value = 594937.4
$ws.Range("I132").Value = 678948.6
$ws.Range("K132").Value = 2036845.8
$ws.Range("M132").Value = -2034315.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 950.2381
$ws.Range("I20").Value = 675.5
$ws.Range("J20").Value = 1829.4
$ws.Range("K20").Value = 675.5
$ws.Range("L20").Value = 1829.4
$ws.Range("M20").Value = -428.5
$ws.Range("N20").Value = -2323.4

$ws.Range("H105").Value = 1607.3334
$ws.Range("I105").Value = 1031.6666
$ws.Range("J105").Value = 2183
$ws.Range("K105").Value = 1031.6666
$ws.Range("L105").Value = 2183
$ws.Range("M105").Value = 715.3334
$ws.Range("N105").Value = -5677

$ws.Range("H134").Value = 781085.4
$ws.Range("I134").Value = 1152467.9
$ws.Range("J134").Value = 293645.94
$ws.Range("K134").Value = 3457403.7
$ws.Range("L134").Value = 880937.8200000001
$ws.Range("M134").Value = -3454868.7
$ws.Range("N134").Value = -886007.8200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1066.6666
$ws.Range("I22").Value = 1001.53845
$ws.Range("J22").Value = 1490
$ws.Range("K22").Value = 1001.53845
$ws.Range("L22").Value = 1490
$ws.Range("M22").Value = -651.53845
$ws.Range("N22").Value = -2190

$ws.Range("H58").Value = 774940.9399999999
$ws.Range("I58").Value = 1544393
$ws.Range("K58").Value = 1544393
$ws.Range("M58").Value = -1544190

$ws.Range("H62").Value = 15000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 15000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 15000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -16248

$ws.Range("H65").Value = 15000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 15000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 75000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -81240

$ws.Range("H122").Value = 3069.8
$ws.Range("I122").Value = 2136.3333
$ws.Range("J122").Value = 3469.8572
$ws.Range("K122").Value = 6408.999899999999
$ws.Range("L122").Value = 10409.5716
$ws.Range("M122").Value = -3958.999899999999
$ws.Range("N122").Value = -15309.5716

$ws.Range("H136").Value = 774940.9399999999
$ws.Range("I136").Value = 1544393
$ws.Range("K136").Value = 4633179
$ws.Range("M136").Value = -4630629

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 20666.834
$ws.Range("J9").Value = 16800
$ws.Range("L9").Value = 50400
$ws.Range("N9").Value = -50848

$ws.Range("H129").Value = 1339.8823
$ws.Range("I129").Value = 562.2143
$ws.Range("J129").Value = 4969
$ws.Range("K129").Value = 1686.6429
$ws.Range("L129").Value = 14907
$ws.Range("M129").Value = 3313.3571
$ws.Range("N129").Value = -24907

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6635.2856
$ws.Range("I70").Value = 6574.5
$ws.Range("K70").Value = 6574.5
$ws.Range("M70").Value = -6304.5

$ws.Range("H73").Value = 6635.2856
$ws.Range("I73").Value = 6574.5
$ws.Range("K73").Value = 6574.5
$ws.Range("M73").Value = -5638.5

$ws.Range("H80").Value = 2181.5134
$ws.Range("I80").Value = 2240.1428
$ws.Range("J80").Value = 2104.5625
$ws.Range("K80").Value = 2240.1428
$ws.Range("L80").Value = 2104.5625
$ws.Range("M80").Value = -1242.1428
$ws.Range("N80").Value = -4100.5625

$ws.Range("H83").Value = 2181.5134
$ws.Range("I83").Value = 2240.1428
$ws.Range("J83").Value = 2104.5625
$ws.Range("K83").Value = 11200.714
$ws.Range("L83").Value = 10522.8125
$ws.Range("M83").Value = -6208.714
$ws.Range("N83").Value = -20506.8125

$ws.Range("H99").Value = 11431.714
$ws.Range("I99").Value = 9828.333000000001
$ws.Range("J99").Value = 21052
$ws.Range("K99").Value = 9828.333000000001
$ws.Range("L99").Value = 21052
$ws.Range("M99").Value = -7582.333000000001
$ws.Range("N99").Value = -25544

$ws.Range("H122").Value = 31999.568
$ws.Range("I122").Value = 44285.125
$ws.Range("J122").Value = 9318.538
$ws.Range("K122").Value = 132855.375
$ws.Range("L122").Value = 27955.614
$ws.Range("M122").Value = -130405.375
$ws.Range("N122").Value = -32855.614

$ws.Range("H126").Value = 928971.2
$ws.Range("I126").Value = 1854262.9
$ws.Range("K126").Value = 5562788.699999999
$ws.Range("M126").Value = -5560318.699999999

$ws.Range("H132").Value = 194098.34
$ws.Range("I132").Value = 241352.44
$ws.Range("K132").Value = 724057.3200000001
$ws.Range("M132").Value = -721527.3200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3916.7896
$ws.Range("I40").Value = 3776.1875
$ws.Range("K40").Value = 3776.1875
$ws.Range("M40").Value = -3640.1875

$ws.Range("H122").Value = 5207.5654
$ws.Range("I122").Value = 5001.316
$ws.Range("J122").Value = 6187.25
$ws.Range("K122").Value = 15003.948
$ws.Range("L122").Value = 18561.75
$ws.Range("M122").Value = -12553.948
$ws.Range("N122").Value = -23461.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1416.069
$ws.Range("I113").Value = 270.8
$ws.Range("J113").Value = 2643.1428
$ws.Range("K113").Value = 812.4000000000001
$ws.Range("L113").Value = 7929.428400000001
$ws.Range("M113").Value = 1357.6
$ws.Range("N113").Value = -12269.4284

$ws.Range("H122").Value = 2421.157
$ws.Range("I122").Value = 2100.9565
$ws.Range("J122").Value = 5367
$ws.Range("K122").Value = 6302.869499999999
$ws.Range("L122").Value = 16101
$ws.Range("M122").Value = -3852.869499999999
$ws.Range("N122").Value = -21001

$ws.Range("H132").Value = 9121652
$ws.Range("I132").Value = 10845818
$ws.Range("K132").Value = 32537454
$ws.Range("M132").Value = -32534924

$ws.Range("H136").Value = 28548696
$ws.Range("I136").Value = 32424638
$ws.Range("J136").Value = 125114.664
$ws.Range("K136").Value = 97273914
$ws.Range("L136").Value = 375343.992
$ws.Range("M136").Value = -97271364
$ws.Range("N136").Value = -380443.992
